$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# DCS approval d_table: update the FFI Employee ID value in A2
# from "fff1616" to "FFI1311"
$ws.Range("A2").Value = "FFI1311"
